$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) AMI $ amount columns (L, M, N, O) were stored as text like "$85,000".
#    Fix: make them real numbers, formatted with thousands separators.
# ------------------------------------------------------------------
$ws.Range("L2").Value = 85000
$ws.Range("L2").NumberFormat = "#,##0"

$ws.Range("M2").Value = 24000
$ws.Range("M2").NumberFormat = "#,##0"

$ws.Range("N2").Value = 45000

$ws.Range("O2").Value = 67000

# ------------------------------------------------------------------
# 2) Insert a new "vacancy_rate" column before the old "vacancy_rates"
#    column (currently column W), so the old text-based "6.5%" column
#    is replaced by a true numeric percentage column.
# ------------------------------------------------------------------
$ws.Columns.Item("W").EntireColumn.Insert()

$ws.Range("W1").Value = "vacancy_rate"
$ws.Range("W2").Value = 0.065
$ws.Range("W2").NumberFormat = "0.00%"

# Remove the old duplicate "vacancy_rates" column, now shifted to X.
$ws.Columns.Item("X").EntireColumn.Delete()

# ------------------------------------------------------------------
# 3) Column width / view tweaks observed in the saved file.
# ------------------------------------------------------------------
$ws.Columns.Item("M").ColumnWidth = 24
$ws.Columns.Item("W").ColumnWidth = 15.21875

$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("W2").Select()
